# -----------------------------------------------------------------------
# C5-PowerPoint.pptx edit
#
# 1) The table on slide 6 (the "Source / ... " table) gets a new built-in
#    table style applied (GUID changes from the "Medium Style 2 - Accent 1"
#    id to the new id used by the author).
#
# 2) The presentation's design/theme colour scheme is swapped from the
#    "Integral" palette to the plain "Office" palette (the two theme parts
#    in the package effectively trade their colour schemes). We reproduce
#    this by writing each of the twelve theme colour slots (dk1, lt1, dk2,
#    lt2, accent1-6, hlink, folHlink) on the presentation's theme colour
#    scheme to the standard "Office" RGB values, matching the new
#    ppt/theme/theme2.xml contents bit for bit. (Font scheme / format
#    scheme are untouched because they are already identical between the
#    two themes.)
# -----------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Table style -----------------------------------------------------
$tableSlide = $p.Slides.Item(6)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{E16F3BFD-A8A4-43C6-828A-D433A52623C9}")

# --- 2) Theme colours -----------------------------------------------------
# VBA/PowerPoint RGB() encodes colours as R + G*256 + B*65536, i.e. the
# reverse byte order of the "RRGGBB" hex string stored in the OOXML
# <a:srgbClr val="RRGGBB"/> element.
function HexToVbaRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the standard Office theme colours.
$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$themeColorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $themeColorScheme.Count; $i++) {
    $themeColorScheme.Colors($i).RGB = HexToVbaRgb($officeColors[$i - 1])
}
